$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 becomes the "radio/dropdown" header pair: No | Mango
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Mango"

# Give B1 the same direct formatting as A1 (font/alignment) without touching
# the value we just wrote, so both header cells share one style.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# The old sample rows (test123 / second test / 3rd test / test 5) go away;
# only A1/B1 keep data now, so fully clear out the leftover cells.
$ws.Range("A2").Clear()
$ws.Range("A3").Clear()
$ws.Range("A5").Clear()
